$p = $ppt.ActivePresentation

# Add a new slide at the end (index 4) using the "Title and Content" layout,
# matching the layout used by the other topic slides (slideLayout2.xml).
$s = $p.Slides.Add(4, 2)

# Title
$s.Shapes.Item(1).TextFrame.TextRange.InsertAfter("Benefícios e Desafios da Inteligência Artificia") | Out-Null

# Body placeholder with multiple outline levels
$tf = $s.Shapes.Item(2).TextFrame
$tf.TextRange.InsertAfter("Listagem dos principais benefícios da IA:`rAutomação de tarefas`rMelhoria da eficiência`rPrevisões precisas`rDestaque para os desafios:`rPrivacidade`rÉtica`rImpactos no mercado de trabalho") | Out-Null

$tf.TextRange.Paragraphs(2).IndentLevel = 2
$tf.TextRange.Paragraphs(3).IndentLevel = 2
$tf.TextRange.Paragraphs(4).IndentLevel = 2
$tf.TextRange.Paragraphs(6).IndentLevel = 2
$tf.TextRange.Paragraphs(7).IndentLevel = 2
$tf.TextRange.Paragraphs(8).IndentLevel = 2
